$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row new values for columns HK (col 219, existing - gains explicit style),
# HL (220, new), HM (221, new), HN (222, new, no explicit style - matches
# the "last column has no explicit style" pattern used throughout the sheet).
$data = @(
    @(10240, 10248, 10258, 10265),
    @(2020, 2020, 2020, 2020),
    @(9, 10, 11, 12),
    @(1, 0, 0, 1),
    @(0, 1, 1, 0),
    @(50, 37, 38, 54),
    @(119, 88, 62, 111),
    @(-69, -51, -24, -57),
    @(0, 0, 0, 0),
    @(12, 11, 4, 18),
    @(137, 133, 158, 133),
    @(121, 119, 119, 119),
    @(258, 252, 277, 252),
    @(1.13, 1.12, 1.33, 1.12),
    @(53, 38, 72, 50),
    @(73, 67, 52, 40),
    @(27, 31, 27, 45),
    @(22, 25, 12, 17),
    @(17, 17, 20, 18),
    @(7, 5, 5, 8),
    @(3, 2, 3, 6),
    @(3, 4, 5, 6),
    @(5, 3, 3, 0),
    @(15, 12, 13, 14),
    @(46.7, 41.7, 38.5, 57.1),
    @(36.86, 50.4, 55.4, 31.5),
    @(17.2, 21, 21.31, 18),
    @(20, 28, 29, 19),
    @(47, 44, 54, 39),
    @(29, 31, 31, 43),
    @(36, 33, 30, 32),
    @(2.4, 2.75, 2.31, 2.29),
    @(5.14, 6.6, 6, 4),
    @(27.8, 27.3, 33.3, 43.8),
    @(19.4, 15.2, 16.7, 25),
    @(187.2, 187.3, 187, 186.5),
    @(85.3, 85.2, 85.5, 85.9),
    @(25.33, 24.33, 24.74, 25.24),
    @(83, 74.1, 73.1, 87.7),
    @(10, 13, 13, 11),
    @(2, 1, 0, 1),
    @(5, 3, 3, 3),
    @(5, 5, 6, 7),
    @(102, 115, 113, 99),
    @(155, 123, 164, 146),
    @(185, 168, 198, 180),
    @(71.7, 66.7, 71.5, 71.4),
    @(47, 44, 54, 39),
    @(5, 11, 13, 5),
    @(5, 4, 5, 3),
    @(20, 28, 29, 19),
    @(29, 31, 31, 43),
    @(38, 36, 36, 32),
    @(5, 4, 1, 2),
    @(3, 2, 3, 6),
    @(42.9, 40, 60, 75),
    @(201, 176, 185, 201),
    @(157, 114, 139, 172),
    @(358, 290, 324, 373),
    @(1.28, 1.54, 1.33, 1.17),
    @(100, 66, 77, 78),
    @(66, 56, 60, 54),
    @(29, 43, 28, 20),
    @(17, 17, 20, 18),
    @(22, 25, 12, 17),
    @(19, 13, 10, 16),
    @(13, 6, 5, 11),
    @(3, 8, 2, 14),
    @(2, 2, 0, 1),
    @(24, 23, 12, 31),
    @(79.2, 56.5, 83.3, 51.6),
    @(18.84, 22.31, 32.4, 23.31),
    @(14.92, 12.61, 27, 12.03),
    @(35, 35, 32, 44),
    @(53, 54, 48, 40),
    @(28, 28, 25, 24),
    @(48, 45, 42, 61),
    @(2, 1.96, 3.5, 1.97),
    @(2.53, 3.46, 4.2, 3.81),
    @(45.8, 46.7, 28.6, 49.2),
    @(39.6, 28.9, 23.8, 26.2),
    @(188.9, 187.5, 188.2, 186.2),
    @(87.2, 87.6, 87.6, 84),
    @(25.33, 25.66, 26.24, 24.58),
    @(86.2, 96.1, 94, 74.3),
    @(9, 5, 10, 9),
    @(5, 9, 3, 6),
    @(2, 5, 4, 5),
    @(6, 3, 5, 2),
    @(125, 135, 108, 138),
    @(232, 147, 214, 236),
    @(269, 192, 242, 280),
    @(75.1, 66.2, 74.7, 75.1),
    @(53, 54, 48, 40),
    @(15, 15, 9, 10),
    @(10, 10, 11, 17),
    @(35, 35, 32, 44),
    @(28, 28, 25, 24),
    @(33, 39, 42, 36),
    @(7, 2, 2, 18),
    @(13, 6, 5, 11),
    @(68.4, 46.2, 50, 68.8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $vals = $data[$i]

    # HK already has a value; clear then re-set so the engine re-applies
    # the column's default style (s="1") explicitly, matching the target.
    $ws.Cells.Item($r, 219).ClearContents()
    $ws.Cells.Item($r, 219).Value = $vals[0]

    # HL and HM are brand new cells; writing a value picks up the column's
    # default style (s="1") automatically.
    $ws.Cells.Item($r, 220).Value = $vals[1]
    $ws.Cells.Item($r, 221).Value = $vals[2]

    # HN is also new, but must NOT carry an explicit style attribute (it is
    # the new "last" column, mirroring how HK looked before the edit).
    # Resetting alignment to the default General/Bottom collapses the style
    # back to the sheet's implicit default (no s attribute written).
    $ws.Cells.Item($r, 222).Value = $vals[3]
    $ws.Cells.Item($r, 222).HorizontalAlignment = -4128
    $ws.Cells.Item($r, 222).VerticalAlignment = -4107
}
